$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") entirely; subsequent rows (3..63) shift up to (2..62).
$ws.Rows.Item(2).Delete()
